$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.971.60'
$ws.Range("E2").Value = '  +2.54%  '
$ws.Range("D3").Value = '1.887.72'
$ws.Range("E3").Value = '  +2.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9968'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.88'
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9966'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4791'
$ws.Range("E7").Value = '  +2.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2847'
$ws.Range("E8").Value = '  +5.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06539'
$ws.Range("E9").Value = '  +4.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.78'
$ws.Range("E10").Value = '  +17.03%  '
$ws.Range("D11").Value = '1.858.79'
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '95.91'
$ws.Range("E12").Value = '  +14.64%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07524'
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.121'
$ws.Range("E14").Value = '  +3.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6543'
$ws.Range("E15").Value = '  +5.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '297.19'
$ws.Range("E16").Value = '  +31.10%  '
$ws.Range("D17").Value = '30.864.09'
$ws.Range("E17").Value = '  +2.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.17'
$ws.Range("E18").Value = '  +6.90%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007569'
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("D21").Value = '2.117.21'
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9968'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.188'
$ws.Range("E23").Value = '  +6.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.141'
$ws.Range("E24").Value = '  +5.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '168.92'
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.285'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.57'
$ws.Range("E27").Value = '  +9.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.971'
$ws.Range("E28").Value = '  +5.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1054'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.374'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.147'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.969'
$ws.Range("E32").Value = '  +4.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04995'
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.181'
$ws.Range("E34").Value = '  +4.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7249'
$ws.Range("E35").Value = '  +2.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.699'
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01942'
$ws.Range("E37").Value = '  +2.91%  '
$ws.Range("E38").Value = '  +2.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.066'
$ws.Range("E39").Value = '  +7.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8972'
$ws.Range("E40").Value = '  +0.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '107.55'
$ws.Range("E41").Value = '  +3.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9976'
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4220'
$ws.Range("E43").Value = '  +5.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.591'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.385'
$ws.Range("E45").Value = '  +5.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.65'
$ws.Range("E46").Value = '  +9.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1236'
$ws.Range("E47").Value = '  +3.62%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.903'
$ws.Range("E48").Value = '  +4.79%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.67'
$ws.Range("E49").Value = '  +5.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.395'
$ws.Range("E50").Value = '  +3.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05573'
$ws.Range("E51").Value = '  +1.24%  '
